$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.911.80'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.312.13'
$ws.Range("E3").Value = '  +3.17%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '96.12'
$ws.Range("E5").Value = '  +1.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.21'
$ws.Range("E6").Value = '  -0.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +0.94%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  -4.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.74'
$ws.Range("E10").Value = '  -6.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0949'
$ws.Range("E11").Value = '  -0.82%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.91'
$ws.Range("E12").Value = '  -6.65%  '

$ws.Range("E13").Value = '  +0.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.668.25'
$ws.Range("E14").Value = '  +3.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.38'
$ws.Range("E15").Value = '  -0.23%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.864'
$ws.Range("E16").Value = '  +4.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.320.56'
$ws.Range("E17").Value = '  +3.78%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.892.52'
$ws.Range("E18").Value = '  +0.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000108'
$ws.Range("E19").Value = '  +2.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.33'
$ws.Range("E20").Value = '  +1.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.56'
$ws.Range("E21").Value = '  +3.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.34'
$ws.Range("E22").Value = '  +0.76%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.36'
$ws.Range("E23").Value = '  +1.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.21'
$ws.Range("E24").Value = '  -1.36%  '

$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.52'
$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.38'
$ws.Range("E27").Value = '  -3.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.48'
$ws.Range("E28").Value = '  -1.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.29'
$ws.Range("E29").Value = '  +1.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.25'
$ws.Range("E30").Value = '  -6.07%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.82'
$ws.Range("E31").Value = '  +1.95%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.16'
$ws.Range("E32").Value = '  +5.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0896'
$ws.Range("E33").Value = '  -4.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.50'
$ws.Range("E34").Value = '  -0.43%  '

$ws.Range("E35").Value = '  +1.16%  '

$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.110'
$ws.Range("E36").Value = '  -3.99%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.46'
$ws.Range("E37").Value = '  +3.07%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0356'
$ws.Range("E38").Value = '  +0.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.32'
$ws.Range("E39").Value = '  -7.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.38'
$ws.Range("E40").Value = '  +9.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.238'
$ws.Range("E41").Value = '  +2.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.38'
$ws.Range("E42").Value = '  +18.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.21'
$ws.Range("E43").Value = '  -6.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.16'
$ws.Range("E44").Value = '  -0.95%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.09'
$ws.Range("E45").Value = '  +7.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.32'
$ws.Range("E46").Value = '  -2.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.103'
$ws.Range("E47").Value = '  +2.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '100.90'
$ws.Range("E48").Value = '  -1.04%  '

$ws.Range("E49").Value = '  +0.53%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.545.79'
$ws.Range("E50").Value = '  +3.41%  '

$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.187'
$ws.Range("E51").Value = '  +12.50%  '

